$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (avoids Excel auto-coercing numeric-looking
# strings like "1.000" or "0.9998" into actual numbers), while leaving the
# cell's style index unchanged (no lingering "@" number format).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "30.504.25"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.910.55"
$ws.Range("E3").Value = "  -0.52%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.9998"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "244.03"
$ws.Range("E5").Value = "  -1.54%  "

# Row 6 - USDC
Set-TextValue "D6" "1.0000"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4832"
$ws.Range("E7").Value = "  +1.85%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.2886"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.06813"

# Row 10 - Litecoin
Set-TextValue "D10" "111.35"
$ws.Range("E10").Value = "  +5.62%  "

# Row 11 - Solana
Set-TextValue "D11" "19.38"
$ws.Range("E11").Value = "  +5.25%  "

# Row 12 - WrappedEther
Set-TextValue "D12" "1.913.75"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13 - TRON
Set-TextValue "D13" "0.07560"
$ws.Range("E13").Value = "  -1.75%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.376"
$ws.Range("E14").Value = "  +1.43%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.6683"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16 - BitcoinCash
Set-TextValue "D16" "293.54"
$ws.Range("E16").Value = "  +0.41%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "30.503.24"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18 - Avalanche
Set-TextValue "D18" "12.99"

# Row 19 - Dai
Set-TextValue "D19" "1.000"
$ws.Range("E19").Value = "  +0.09%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.000007568"
$ws.Range("E20").Value = "  -0.44%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue "D21" "2.162.81"
$ws.Range("E21").Value = "  -0.25%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.491"
$ws.Range("E22").Value = "  -1.18%  "

# Row 23 - BinanceUSD
Set-TextValue "D23" "0.9990"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24 - Chainlink
Set-TextValue "D24" "6.398"
$ws.Range("E24").Value = "  +0.19%  "

# Row 25 - Cosmos
Set-TextValue "D25" "9.447"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26 - Monero
Set-TextValue "D26" "165.10"
$ws.Range("E26").Value = "  -1.62%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "20.28"
$ws.Range("E27").Value = "  -4.15%  "

# Row 28 - LidoDAOToken
Set-TextValue "D28" "2.074"
$ws.Range("E28").Value = "  -1.67%  "

# Row 29 - Stellar
Set-TextValue "D29" "0.1064"
$ws.Range("E29").Value = "  -0.60%  "

# Row 30 - Toncoin
Set-TextValue "D30" "1.428"
$ws.Range("E30").Value = "  +2.35%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "4.130"
$ws.Range("E31").Value = "  -1.15%  "

# Row 32 - Filecoin
Set-TextValue "D32" "4.051"
$ws.Range("E32").Value = "  -0.41%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.04972"
$ws.Range("E33").Value = "  -1.18%  "

# Row 34 - ImmutableX
Set-TextValue "D34" "0.7347"
$ws.Range("E34").Value = "  -0.70%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.132"

# Row 36 - Frax
$ws.Range("E36").Value = "  +0.12%  "

# Rows 37 and 38 swap content: HuobiToken <-> VeChain
Set-TextValue "B37" "VeChain"
Set-TextValue "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.02036"
$ws.Range("E37").Value = "  -2.08%  "

Set-TextValue "B38" "HuobiToken"
Set-TextValue "C38" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D38" "2.718"
$ws.Range("E38").Value = "  -1.07%  "

# Row 39 - MXToken
Set-TextValue "D39" "2.682"
$ws.Range("E39").Value = "  -0.15%  "

# Row 40 - RenderToken
Set-TextValue "D40" "2.008"
$ws.Range("E40").Value = "  -2.52%  "

# Row 41 - Quant
Set-TextValue "D41" "109.15"
$ws.Range("E41").Value = "  -1.57%  "

# Row 42 - TheSandbox
Set-TextValue "D42" "0.4420"
$ws.Range("E42").Value = "  +0.89%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "0.8632"
$ws.Range("E43").Value = "  -1.71%  "

# Row 44 - FraxShare
Set-TextValue "D44" "5.769"
$ws.Range("E44").Value = "  -2.09%  "

# Row 45 - PaxDollar
Set-TextValue "D45" "0.9997"
$ws.Range("E45").Value = "  -0.02%  "

# Row 46 - Aave
Set-TextValue "D46" "69.06"
$ws.Range("E46").Value = "  +1.89%  "

# Row 47 - Aptos
Set-TextValue "D47" "7.192"
$ws.Range("E47").Value = "  -1.08%  "

# Row 48 - BitcoinSV
Set-TextValue "D48" "48.16"
$ws.Range("E48").Value = "  -0.72%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "9.221"
$ws.Range("E49").Value = "  -1.17%  "

# Row 50 - Algorand
Set-TextValue "D50" "0.1227"
$ws.Range("E50").Value = "  -1.24%  "

# Row 51 - WOONetwork
Set-TextValue "D51" "0.2521"
$ws.Range("E51").Value = "  +1.71%  "
